# ProblemsAndSolution.xlsx update
# - Collapses the old 3-row block (rows 121-123) describing:
#     121: "Print all nodes that are at distance k from a leaf node"
#     122: "ReverseEveryKNodes"
#     123: "GetAllUniqueNodesKDistanceFromRoot"
#   into a 2-row block:
#     121: "UniqueNodes at K Distance From Root,Assuming No Duplicate Values at Nodes"
#          with its solution link moved into column C (Sarath Solution)
#     122: "Reverse Every K Nodes in Linked List"
#          with its solution link moved into column C (Sarath Solution)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the old row 123 (GetAllUniqueNodesKDistanceFromRoot / its gist link) entirely -
# its content is superseded by the consolidated row 121.
$ws.Rows.Item(123).Delete()

# Row 121: new problem title, drop the old Barani-solution cell entirely, add the
# new Sarath-solution link (ideone) in column C.
$ws.Range("A121").Value = "UniqueNodes at K Distance From Root,Assuming No Duplicate Values at Nodes"
$linkStyle = $ws.Range("B121").Style
$ws.Range("B121").Clear()
$ws.Hyperlinks.Add($ws.Range("C121"), "http://ideone.com/ymyLOU")
$ws.Range("C121").Style = $linkStyle

# Row 122: renamed problem title, drop the old Barani-solution cell entirely, add
# the new Sarath-solution link (ideone) in column C.
$ws.Range("A122").Value = "Reverse Every K Nodes in Linked List"
$linkStyle = $ws.Range("B122").Style
$ws.Range("B122").Clear()
$ws.Hyperlinks.Add($ws.Range("C122"), "http://ideone.com/25I7AF")
$ws.Range("C122").Style = $linkStyle

# Restore the selection to the new last-used row beneath the table.
[void]$ws.Range("A123").Select()
